# Data source / layout reshuffle + refreshed figures ("update DB connections
# and arrange all code"): several Item Name / UOM rows were re-arranged and
# the sales figures that go with them were refreshed to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Item Name (column C) re-arrangements ---
$ws.Range("C3").Value  = "Dinafex 120mg Tablet"
$ws.Range("C4").Value  = "Dinafex 180mg Tablet"
$ws.Range("C5").Value  = "Dinafex 60mg Tablet"

$ws.Range("C7").Value  = "Etorix 90mg Tablet"
$ws.Range("C8").Value  = "Etorix 120mg Tablet"
$ws.Range("C9").Value  = "Etorix 60mg Tablet - 40's"

$ws.Range("C11").Value = "Flucloxin 500mg Capsule - 36's"
$ws.Range("C12").Value = "Flucloxin 500mg Capsule"

$ws.Range("C14").Value = "Ketonic 10mg Tablet"
$ws.Range("C15").Value = "Ketonic 30mg Injection"
$ws.Range("C16").Value = "Ketonic 30mg IM/IV Injection - 4's"

$ws.Range("C17").Value = "Kynol D 25mg Tablet"
$ws.Range("C18").Value = "Kynol TR 200mg Capsule"
$ws.Range("C19").Value = "Kynol TR 100mg Capsule"

$ws.Range("C25").Value = "Zithrox 15ml Suspension"
$ws.Range("C26").Value = "Zithrox 250mg Tablet - 6's"
$ws.Range("C27").Value = "Zithrox 500mg Tablet"
$ws.Range("C28").Value = "Zithrox 30ml Dry Suspension"

# --- UOM (column D) re-arrangements, matching the Item Name moves above ---
$ws.Range("D7").Value  = "30's"
$ws.Range("D8").Value  = "20's"
$ws.Range("D9").Value  = "40's"

$ws.Range("D11").Value = "36 's"
$ws.Range("D12").Value = "30 's"

$ws.Range("D14").Value = "20's"
$ws.Range("D15").Value = "5 's"
$ws.Range("D16").Value = "4's"

$ws.Range("D17").Value = "60 's"
$ws.Range("D18").Value = "30 's"
$ws.Range("D19").Value = "50 's"

$ws.Range("D25").Value = "15 ml"
$ws.Range("D26").Value = "6's"
$ws.Range("D27").Value = "6 's"
$ws.Range("D28").Value = "30ml"

# --- Refreshed TP (column BC) values, aligned with the re-ordered items ---
$ws.Range("BC3").Value  = 179.91
$ws.Range("BC4").Value  = 224.89
$ws.Range("BC5").Value  = 78.70999999999999

$ws.Range("BC7").Value  = 269.87
$ws.Range("BC9").Value  = 209.9

$ws.Range("BC11").Value = 284.21
$ws.Range("BC12").Value = 237.74

$ws.Range("BC14").Value = 150.38
$ws.Range("BC15").Value = 206.77
$ws.Range("BC16").Value = 165.41

$ws.Range("BC17").Value = 180.45
$ws.Range("BC18").Value = 224.89
$ws.Range("BC19").Value = 262.37

$ws.Range("BC25").Value = 71.95999999999999
$ws.Range("BC26").Value = 89.95999999999999
$ws.Range("BC27").Value = 136.83
$ws.Range("BC28").Value = 97.45

# --- Other refreshed figures (rows 13, 23, 24) ---
$ws.Range("H13").Value  = 238

$ws.Range("H23").Value  = 1117
$ws.Range("L23").Value  = 10

$ws.Range("E24").Value  = 0
$ws.Range("O24").Value  = 173
$ws.Range("U24").Value  = 173
$ws.Range("AV24").Value = 152
$ws.Range("BA24").Value = 138
$ws.Range("BB24").Value = 30239
$ws.Range("BD24").Value = 0
$ws.Range("BE24").Value = 0
